# Sync attendance_reports: normalize the "Recorded By" (column G) value lists
# so the comma-separated list of recorder identities is sorted in strict
# ordinal (byte/code-point) order, e.g. "dnasr281@gmail.com, System"
# becomes "System, dnasr281@gmail.com" because 'S' (83) sorts before 'd' (100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Compares two strings using ordinal (code point) semantics, since the
# built-in comparison/sort operators in this environment are culture-aware
# (case-insensitive alphabetic) and do not reproduce the desired ordering.
# NOTE: uses uniquely-named loop/local variables everywhere to avoid
# clashing with variables of the same name used by callers.
function Compare-Ordinal($cmpA, $cmpB) {
    $cmpLenA = $cmpA.Length
    $cmpLenB = $cmpB.Length
    $cmpMinLen = [Math]::Min($cmpLenA, $cmpLenB)
    for ($cmpIdx = 0; $cmpIdx -lt $cmpMinLen; $cmpIdx++) {
        $cmpCharA = [int][char]$cmpA[$cmpIdx]
        $cmpCharB = [int][char]$cmpB[$cmpIdx]
        if ($cmpCharA -ne $cmpCharB) {
            return $cmpCharA - $cmpCharB
        }
    }
    return $cmpLenA - $cmpLenB
}

# Simple insertion sort using Compare-Ordinal, since Sort-Object uses
# culture-aware comparison by default.
function Sort-Ordinal($sortItems) {
    $sortList = New-Object System.Collections.ArrayList
    foreach ($sortItem in $sortItems) {
        [void]$sortList.Add($sortItem)
    }
    for ($sortOuter = 1; $sortOuter -lt $sortList.Count; $sortOuter++) {
        $sortKey = $sortList[$sortOuter]
        $sortInner = $sortOuter - 1
        while (($sortInner -ge 0) -and ((Compare-Ordinal $sortList[$sortInner] $sortKey) -gt 0)) {
            $sortList[$sortInner + 1] = $sortList[$sortInner]
            $sortInner = $sortInner - 1
        }
        $sortList[$sortInner + 1] = $sortKey
    }
    return $sortList
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Column G = "Recorded By"
$recordedByCol = 7

for ($rowNum = 2; $rowNum -le $lastRow; $rowNum++) {
    $cell = $ws.Cells.Item($rowNum, $recordedByCol)
    $cellVal = $cell.Value2

    if (($cellVal -ne $null) -and ($cellVal -ne "")) {
        $valParts = $cellVal -split ", "
        if ($valParts.Count -gt 1) {
            $sortedParts = Sort-Ordinal $valParts
            $newCellVal = $sortedParts -join ", "
            if ($newCellVal -ne $cellVal) {
                $cell.Value2 = $newCellVal
            }
        }
    }
}
